$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'42.749.93"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "

$cell = $ws.Range("D3")
$cell.Value = "'2.554.60"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  -0.15%  "

$cell = $ws.Range("D5")
$cell.Value = "'303.14"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.70%  "

$cell = $ws.Range("D6")
$cell.Value = "'98.42"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +6.82%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  +0.02%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.547"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.63%  "

$cell = $ws.Range("D10")
$cell.Value = "'36.77"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.77%  "

$cell = $ws.Range("D11")
$cell.Value = "'0.0808"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("E12").Value = "  +8.04%  "

$cell = $ws.Range("D13")
$cell.Value = "'7.67"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.57%  "

$cell = $ws.Range("D14")
$cell.Value = "'2.529.47"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.13%  "

$cell = $ws.Range("D15")
$cell.Value = "'0.880"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "

$cell = $ws.Range("D16")
$cell.Value = "'14.79"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +4.58%  "

$cell = $ws.Range("D17")
$cell.Value = "'42.800.84"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "

$cell = $ws.Range("D18")
$cell.Value = "'13.27"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +5.50%  "

$cell = $ws.Range("D19")
$cell.Value = "'0.0₃0985"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "

$cell = $ws.Range("D20")
$cell.Value = "'6.58"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.09%  "

$cell = $ws.Range("D21")
$cell.Value = "'71.63"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "

$cell = $ws.Range("D22")
$cell.Value = "'254.56"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.33%  "

$cell = $ws.Range("D23")
$cell.Value = "'2.95"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "

$cell = $ws.Range("D24")
$cell.Value = "'2.09"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "

$cell = $ws.Range("D25")
$cell.Value = "'27.55"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -6.51%  "

$ws.Range("E26").Value = "  +0.06%  "

$cell = $ws.Range("D27")
$cell.Value = "'10.06"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "

$cell = $ws.Range("D28")
$cell.Value = "'37.85"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +3.83%  "

$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("E30").Value = "  -0.20%  "

$cell = $ws.Range("D31")
$cell.Value = "'156.01"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +2.57%  "

$cell = $ws.Range("D32")
$cell.Value = "'2.18"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("E33").Value = "  +1.01%  "

$cell = $ws.Range("D34")
$cell.Value = "'0.0807"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("E35").Value = "  -2.55%  "

$cell = $ws.Range("D36")
$cell.Value = "'26.38"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +9.26%  "

$cell = $ws.Range("D37")
$cell.Value = "'18.54"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +11.58%  "

$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("E39").Value = "  -0.32%  "

$cell = $ws.Range("D40")
$cell.Value = "'2.08"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +32.18%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D41")
$cell.Value = "'3.39"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D42")
$cell.Value = "'3.86"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Range("D43")
$cell.Value = "'2.066.56"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D44")
$cell.Value = "'0.0302"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("E45").Value = "  +0.02%  "

$cell = $ws.Range("D46")
$cell.Value = "'87.84"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +3.17%  "

$cell = $ws.Range("D47")
$cell.Value = "'9.20"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +6.21%  "

$cell = $ws.Range("D48")
$cell.Value = "'2.799.67"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "

$cell = $ws.Range("D49")
$cell.Value = "'75.10"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +8.25%  "

$cell = $ws.Range("D50")
$cell.Value = "'103.19"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "

$cell = $ws.Range("D51")
$cell.Value = "'0.190"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.79%  "

